$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting existing rows 108..205 down to 109..206
$ws.Rows("108:108").Insert()

# Populate the newly inserted row 108 with the new data record
$ws.Cells.Item(108, 1).Value = 4
$ws.Cells.Item(108, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(108, 3).Value = "Los Lagos"
$ws.Cells.Item(108, 4).Value = 44589
$ws.Cells.Item(108, 5).Value = 10
$ws.Cells.Item(108, 6).Value = 100112044
$ws.Cells.Item(108, 7).Value = "Perejil"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 120
$ws.Cells.Item(108, 11).Value = 5000
$ws.Cells.Item(108, 12).Value = 5000
$ws.Cells.Item(108, 13).Value = 5000
$ws.Cells.Item(108, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(108, 15).Value = "Región Metropolitana"
$ws.Cells.Item(108, 16).Value = 1667
$ws.Cells.Item(108, 17).Value = 3
$ws.Cells.Item(108, 18).Value = "Hortaliza"
